$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column at column N (14) on the "Repayment schedule" sheet.
# This shifts the previous N/O/P columns (Outstanding / Paid Date / Disbursement)
# one position to the right, becoming O/P/Q.
$wsRepay.Columns("N:N").Insert()

# Give the newly inserted column the same width as the former "Paid Date" column (11 characters).
$wsRepay.Columns("N:N").ColumnWidth = 10.166666666666666

# Move the active sheet/selection from "Transactions" to "Repayment schedule",
# with the active cell at R8.
$wsRepay.Select()
$wsRepay.Range("R8").Select()
